$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 17: HappinessList ---
$ws.Range("A17").Value2 = "HappinessList"

$ws.Range("B17").Value2 = "(n: int)"
$ws.Range("B17").Characters(5, 3).Font.Color = 4697456

$ws.Range("C17").Value2 = "list"
$ws.Range("C9").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("D17").Value2 = "Creates a list which we will use to give the happiness values"

# --- New row 18: HappinessMeter ---
$ws.Range("A18").Value2 = "HappinessMeter"

$ws.Range("D18").Value2 = "Calculates the total happiness for a given solution in dictionary form (dictionary with courses as the keys)"

$ws.Range("C18").Value2 = "int"
$ws.Range("C2").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("B18").Value2 = "(solution: dict)"
$ws.Range("B18").Characters(12, 4).Font.Color = 4697456
$ws.Range("B18").Characters(16, 1).Font.Color = 0

$excel.CutCopyMode = $false

# --- Resize table to include new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D18"))

# --- Widen column D ---
$ws.Columns.Item(4).ColumnWidth = 87.16666666666667

# --- Update selection to match final cursor position ---
[void]$ws.Range("B18").Select()
